$d = $word.ActiveDocument

# Phase 1: replace each old value with a unique placeholder token to avoid
# any chained/cascading replacements when an old value equals another new value.
$d.Content.Find.Execute("2025-07-01 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_0@@", 2) | Out-Null
$d.Content.Find.Execute("936×5=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_1@@", 2) | Out-Null
$d.Content.Find.Execute("388×6=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_2@@", 2) | Out-Null
$d.Content.Find.Execute("347×4=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_3@@", 2) | Out-Null
$d.Content.Find.Execute("634×9=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_4@@", 2) | Out-Null
$d.Content.Find.Execute("957×7=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_5@@", 2) | Out-Null
$d.Content.Find.Execute("739×9=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_6@@", 2) | Out-Null
$d.Content.Find.Execute("546×2=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_7@@", 2) | Out-Null
$d.Content.Find.Execute("915×7=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_8@@", 2) | Out-Null
$d.Content.Find.Execute("224×6=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_9@@", 2) | Out-Null
$d.Content.Find.Execute("947×7=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_10@@", 2) | Out-Null
$d.Content.Find.Execute("836×5=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_11@@", 2) | Out-Null
$d.Content.Find.Execute("784×2=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_12@@", 2) | Out-Null
$d.Content.Find.Execute("703×4=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_13@@", 2) | Out-Null
$d.Content.Find.Execute("367×7=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_14@@", 2) | Out-Null
$d.Content.Find.Execute("229×3=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_15@@", 2) | Out-Null
$d.Content.Find.Execute("483×2=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_16@@", 2) | Out-Null
$d.Content.Find.Execute("884×5=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_17@@", 2) | Out-Null
$d.Content.Find.Execute("263×3=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_18@@", 2) | Out-Null
$d.Content.Find.Execute("428×9=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_19@@", 2) | Out-Null
$d.Content.Find.Execute("413×9=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_20@@", 2) | Out-Null
$d.Content.Find.Execute("686×4=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_21@@", 2) | Out-Null
$d.Content.Find.Execute("575×4=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_22@@", 2) | Out-Null
$d.Content.Find.Execute("697×2=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_23@@", 2) | Out-Null
$d.Content.Find.Execute("382×8=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_24@@", 2) | Out-Null
$d.Content.Find.Execute("567×6=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_25@@", 2) | Out-Null

# Phase 2: replace each placeholder token with its final new value.
$d.Content.Find.Execute("@@PLACEHOLDER_0@@", $true, $false, $false, $false, $false, $true, 1, $false, "2025-07-02 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_1@@", $true, $false, $false, $false, $false, $true, 1, $false, "276×8=", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_2@@", $true, $false, $false, $false, $false, $true, 1, $false, "977×5=", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_3@@", $true, $false, $false, $false, $false, $true, 1, $false, "486×7=", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_4@@", $true, $false, $false, $false, $false, $true, 1, $false, "768×8=", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_5@@", $true, $false, $false, $false, $false, $true, 1, $false, "310×8=", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_6@@", $true, $false, $false, $false, $false, $true, 1, $false, "413×9=", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_7@@", $true, $false, $false, $false, $false, $true, 1, $false, "565×2=", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_8@@", $true, $false, $false, $false, $false, $true, 1, $false, "710×8=", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_9@@", $true, $false, $false, $false, $false, $true, 1, $false, "275×7=", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_10@@", $true, $false, $false, $false, $false, $true, 1, $false, "688×7=", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_11@@", $true, $false, $false, $false, $false, $true, 1, $false, "661×5=", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_12@@", $true, $false, $false, $false, $false, $true, 1, $false, "147×5=", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_13@@", $true, $false, $false, $false, $false, $true, 1, $false, "725×7=", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_14@@", $true, $false, $false, $false, $false, $true, 1, $false, "930×2=", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_15@@", $true, $false, $false, $false, $false, $true, 1, $false, "206×9=", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_16@@", $true, $false, $false, $false, $false, $true, 1, $false, "959×9=", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_17@@", $true, $false, $false, $false, $false, $true, 1, $false, "138×8=", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_18@@", $true, $false, $false, $false, $false, $true, 1, $false, "655×3=", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_19@@", $true, $false, $false, $false, $false, $true, 1, $false, "354×7=", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_20@@", $true, $false, $false, $false, $false, $true, 1, $false, "114×5=", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_21@@", $true, $false, $false, $false, $false, $true, 1, $false, "611×8=", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_22@@", $true, $false, $false, $false, $false, $true, 1, $false, "249×7=", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_23@@", $true, $false, $false, $false, $false, $true, 1, $false, "318×5=", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_24@@", $true, $false, $false, $false, $false, $true, 1, $false, "926×3=", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_25@@", $true, $false, $false, $false, $false, $true, 1, $false, "383×3=", 2) | Out-Null
